# Trade #16 closed at 2026-02-16 21:58:11 - leadlag UP +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1000.09   # Current Capital
$summary.Range("B5").Value = 0.22      # Total P&L %
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (leadlag row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.07
$status.Range("D5").Value = 8
$status.Range("F5").Value = 0.07000000000000001
$status.Range("G5").Value = 37.5

# ---------------------------------------------------------------------
# Helper to update an existing "open trade" row once it is closed
# ---------------------------------------------------------------------
function Update-ClosedTrade($ws, $row) {
    $ws.Cells.Item($row, 7).Value = 68318.23723100001   # Exit Price (G)
    $ws.Cells.Item($row, 8).Value = "CLOSED"             # Status (H)
    $ws.Cells.Item($row, 9).Value = 0.0368               # P&L % (I)
    $ws.Cells.Item($row, 11).Value = 100.07              # Capital After (K)
    $ws.Cells.Item($row, 14).Value = "time_exit_5min"    # Exit Reason (N)
    $ws.Cells.Item($row, 15).Value = 5                   # Duration (min) (O)
}

# Helper to append the newly opened trade #16
function Add-NewTrade($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 16
    # Leading apostrophe forces Excel to keep the text "2026-02-16" instead
    # of auto-converting it to a date serial; reset the style afterwards so
    # the cell is not left with a "quote prefix" format flag.
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "21:58:11"
    $ws.Cells.Item($row, 4).Value = "leadlag"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 68383.37
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.0660986376279
    $ws.Cells.Item($row, 12).Value = 0.6758999999999999
    $ws.Cells.Item($row, 13).Value = "Coinbase leading with 0.068% move"
    # Exit Price (G) / Exit Reason (N) stay blank until the trade closes;
    # still materialize the cells so the row matches the sheet's shape.
    $ws.Cells.Item($row, 7).NumberFormat = "General"
    $ws.Cells.Item($row, 7).Style = "Normal"
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 14).NumberFormat = "General"
    $ws.Cells.Item($row, 14).Style = "Normal"
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Update-ClosedTrade $allTrades 10
Add-NewTrade $allTrades 17

# ---------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")
Update-ClosedTrade $leadlag 9
Add-NewTrade $leadlag 16

Write-Host "Applied trade #16 update"
